# Remove the "- Confidence" columns (B, D, F) from every worksheet in the
# workbook, leaving just File Number / Dated / Effective (and their values)
# in columns A/B/C.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Delete from right to left so earlier deletions don't shift the
    # indices of columns we still need to remove.
    $ws.Columns.Item(6).Delete()
    $ws.Columns.Item(4).Delete()
    $ws.Columns.Item(2).Delete()
}
